$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E6").Value = "On Hold"
$ws.Range("E28").Value = "On Hold (I-SPY2 data not yet available)"
$ws.Range("E30").Value = "Complete"
$ws.Range("E31").Value = "Complete"

$ws.Range("A32").Value = 31
$ws.Range("B32").Value = "Send the Agilent 415K ADF custom array design to Zhong."
$ws.Range("C32").Value = "Henry Schaefer"
$ws.Range("D32").Value = 39491
$ws.Range("E32").Value = "Not Started"

$ws.Range("B38").Select()
